# Auto-generated edit script: updates cached currentAveragePrice / LevePrice /
# LeveProfit values across all 8 Leve-profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 178.625
$ws.Range("I5").Value = 109.666664
$ws.Range("J5").Value = 220
$ws.Range("K5").Value = 109.666664
$ws.Range("L5").Value = 220
$ws.Range("M5").Value = 5.333336000000003
$ws.Range("N5").Value = -450
$ws.Range("H98").Value = 634.9375
$ws.Range("I98").Value = 650.2
$ws.Range("J98").Value = 406
$ws.Range("K98").Value = 650.2
$ws.Range("L98").Value = 406
$ws.Range("M98").Value = 847.8
$ws.Range("N98").Value = -3402
$ws.Range("H111").Value = 1027.9
$ws.Range("I111").Value = 566.3333
$ws.Range("J111").Value = 1225.7142
$ws.Range("K111").Value = 1698.9999
$ws.Range("L111").Value = 3677.1426
$ws.Range("M111").Value = 1368.0001
$ws.Range("N111").Value = -9811.142599999999
$ws.Range("H113").Value = 3287.353
$ws.Range("I113").Value = 3318.6365
$ws.Range("J113").Value = 3230
$ws.Range("K113").Value = 3318.6365
$ws.Range("L113").Value = 3230
$ws.Range("M113").Value = -64.63650000000007
$ws.Range("N113").Value = -9738
$ws.Range("H116").Value = 1818.7273
$ws.Range("I116").Value = 1625
$ws.Range("J116").Value = 2051.2
$ws.Range("K116").Value = 1625
$ws.Range("L116").Value = 2051.2
$ws.Range("M116").Value = 1817
$ws.Range("N116").Value = -8935.200000000001
$ws.Range("H122").Value = 634.9375
$ws.Range("I122").Value = 650.2
$ws.Range("J122").Value = 406
$ws.Range("K122").Value = 1950.6
$ws.Range("L122").Value = 1218
$ws.Range("M122").Value = 499.3999999999999
$ws.Range("N122").Value = -6118
$ws.Range("H129").Value = 934.4138
$ws.Range("I129").Value = 752.9375
$ws.Range("J129").Value = 1003.5476
$ws.Range("K129").Value = 2258.8125
$ws.Range("L129").Value = 3010.6428
$ws.Range("M129").Value = 2741.1875
$ws.Range("N129").Value = -13010.6428
$ws.Range("H132").Value = 9806926
$ws.Range("I132").Value = 10377.333
$ws.Range("J132").Value = 24501750
$ws.Range("K132").Value = 31131.999
$ws.Range("L132").Value = 73505250
$ws.Range("M132").Value = -28601.999
$ws.Range("N132").Value = -73510310

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1986.9584
$ws.Range("I2").Value = 1648.3636
$ws.Range("K2").Value = 1648.3636
$ws.Range("M2").Value = -1535.3636
$ws.Range("H32").Value = 16333.88
$ws.Range("I32").Value = 16516.297
$ws.Range("J32").Value = 14834
$ws.Range("K32").Value = 16516.297
$ws.Range("L32").Value = 14834
$ws.Range("M32").Value = -16229.297
$ws.Range("N32").Value = -15408
$ws.Range("H46").Value = 2181.3333
$ws.Range("I46").Value = 1949.5
$ws.Range("K46").Value = 1949.5
$ws.Range("M46").Value = -1630.5
$ws.Range("H61").Value = 71573750
$ws.Range("I61").Value = 143001760
$ws.Range("J61").Value = 145754.28
$ws.Range("K61").Value = 143001760
$ws.Range("L61").Value = 145754.28
$ws.Range("M61").Value = -143001548
$ws.Range("N61").Value = -146178.28
$ws.Range("H74").Value = 10082848
$ws.Range("I74").Value = 14767205
$ws.Range("J74").Value = 128587.5
$ws.Range("K74").Value = 14767205
$ws.Range("L74").Value = 128587.5
$ws.Range("M74").Value = -14766331
$ws.Range("N74").Value = -130335.5
$ws.Range("H77").Value = 10082848
$ws.Range("I77").Value = 14767205
$ws.Range("J77").Value = 128587.5
$ws.Range("K77").Value = 73836025
$ws.Range("L77").Value = 642937.5
$ws.Range("M77").Value = -73831657
$ws.Range("N77").Value = -651673.5
$ws.Range("H116").Value = 1986.9584
$ws.Range("I116").Value = 1648.3636
$ws.Range("K116").Value = 1648.3636
$ws.Range("M116").Value = 645.6364000000001
$ws.Range("H122").Value = 11113238
$ws.Range("I122").Value = 2516.6667
$ws.Range("J122").Value = 27779320
$ws.Range("K122").Value = 7550.000100000001
$ws.Range("L122").Value = 83337960
$ws.Range("M122").Value = -5100.000100000001
$ws.Range("N122").Value = -83342860
$ws.Range("H132").Value = 101331.91
$ws.Range("I132").Value = 77742.92
$ws.Range("J132").Value = 139664
$ws.Range("K132").Value = 233228.76
$ws.Range("L132").Value = 418992
$ws.Range("M132").Value = -230698.76
$ws.Range("N132").Value = -424052
$ws.Range("H136").Value = 71573750
$ws.Range("I136").Value = 143001760
$ws.Range("J136").Value = 145754.28
$ws.Range("K136").Value = 429005280
$ws.Range("L136").Value = 437262.84
$ws.Range("M136").Value = -429002730
$ws.Range("N136").Value = -442362.84

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1986.9584
$ws.Range("I3").Value = 1648.3636
$ws.Range("K3").Value = 1648.3636
$ws.Range("M3").Value = -1534.3636
$ws.Range("H134").Value = 3350.6538
$ws.Range("I134").Value = 2616.9473
$ws.Range("J134").Value = 5342.143
$ws.Range("K134").Value = 7850.841899999999
$ws.Range("L134").Value = 16026.429
$ws.Range("M134").Value = -5315.841899999999
$ws.Range("N134").Value = -21096.429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3026.4707
$ws.Range("I99").Value = 2683.75
$ws.Range("J99").Value = 3331.111
$ws.Range("K99").Value = 2683.75
$ws.Range("L99").Value = 3331.111
$ws.Range("M99").Value = -1185.75
$ws.Range("N99").Value = -6327.111
$ws.Range("H107").Value = 423.94446
$ws.Range("J107").Value = 450.625
$ws.Range("L107").Value = 450.625
$ws.Range("N107").Value = -4290.625
$ws.Range("H122").Value = 1658.2759
$ws.Range("I122").Value = 1266.238
$ws.Range("J122").Value = 2687.375
$ws.Range("K122").Value = 3798.714
$ws.Range("L122").Value = 8062.125
$ws.Range("M122").Value = -1348.714
$ws.Range("N122").Value = -12962.125
$ws.Range("H126").Value = 3026.4707
$ws.Range("I126").Value = 2683.75
$ws.Range("J126").Value = 3331.111
$ws.Range("K126").Value = 8051.25
$ws.Range("L126").Value = 9993.332999999999
$ws.Range("M126").Value = -5581.25
$ws.Range("N126").Value = -14933.333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 800.2174
$ws.Range("I5").Value = 700
$ws.Range("J5").Value = 821.3158
$ws.Range("K5").Value = 2100
$ws.Range("L5").Value = 2463.9474
$ws.Range("M5").Value = -1988
$ws.Range("N5").Value = -2687.9474
$ws.Range("H12").Value = 38462124
$ws.Range("J12").Value = 262.42856
$ws.Range("L12").Value = 787.28568
$ws.Range("N12").Value = -1133.28568
$ws.Range("H14").Value = 894.88
$ws.Range("I14").Value = 894.88
$ws.Range("K14").Value = 2684.64
$ws.Range("M14").Value = -2511.64
$ws.Range("H68").Value = 1191.4828
$ws.Range("J68").Value = 1230.2
$ws.Range("L68").Value = 3690.6
$ws.Range("N68").Value = -5312.6
$ws.Range("H71").Value = 1191.4828
$ws.Range("J71").Value = 1230.2
$ws.Range("L71").Value = 11071.8
$ws.Range("N71").Value = -19183.8
$ws.Range("H107").Value = 638.5238000000001
$ws.Range("I107").Value = 831
$ws.Range("J107").Value = 463.54544
$ws.Range("K107").Value = 2493
$ws.Range("L107").Value = 1390.63632
$ws.Range("M107").Value = -573
$ws.Range("N107").Value = -5230.63632
$ws.Range("H122").Value = 848
$ws.Range("I122").Value = 449.08334
$ws.Range("J122").Value = 1645.8334
$ws.Range("K122").Value = 4041.75006
$ws.Range("L122").Value = 14812.5006
$ws.Range("M122").Value = -1591.75006
$ws.Range("N122").Value = -19712.5006
$ws.Range("H131").Value = 940.60565
$ws.Range("J131").Value = 1007.93445
$ws.Range("L131").Value = 3023.80335
$ws.Range("N131").Value = -13103.80335
$ws.Range("H135").Value = 800.2174
$ws.Range("I135").Value = 700
$ws.Range("J135").Value = 821.3158
$ws.Range("K135").Value = 6300
$ws.Range("L135").Value = 7391.8422
$ws.Range("M135").Value = -3765
$ws.Range("N135").Value = -12461.8422

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2499.6667
$ws.Range("I102").Value = 1799.6
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 1799.6
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -177.5999999999999
$ws.Range("N102").Value = -9244
$ws.Range("H113").Value = 1139.3889
$ws.Range("I113").Value = 923.5
$ws.Range("J113").Value = 1895
$ws.Range("K113").Value = 923.5
$ws.Range("L113").Value = 1895
$ws.Range("M113").Value = 1246.5
$ws.Range("N113").Value = -6235
$ws.Range("H122").Value = 2937.875
$ws.Range("I122").Value = 1997.5
$ws.Range("J122").Value = 3251.3333
$ws.Range("K122").Value = 5992.5
$ws.Range("L122").Value = 9753.999899999999
$ws.Range("M122").Value = -3542.5
$ws.Range("N122").Value = -14653.9999
$ws.Range("H126").Value = 1872.1111
$ws.Range("I126").Value = 1533.5714
$ws.Range("K126").Value = 4600.7142
$ws.Range("M126").Value = -2130.7142
$ws.Range("H132").Value = 102392.7
$ws.Range("I132").Value = 73283.28999999999
$ws.Range("K132").Value = 219849.87
$ws.Range("M132").Value = -217319.87

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1825
$ws.Range("I7").Value = 1828.5714
$ws.Range("J7").Value = 1800
$ws.Range("K7").Value = 1828.5714
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = -1716.5714
$ws.Range("N7").Value = -2024
$ws.Range("H61").Value = 2279.5
$ws.Range("I61").Value = 2069.9
$ws.Range("K61").Value = 2069.9
$ws.Range("M61").Value = -1867.9
$ws.Range("H113").Value = 2279.5
$ws.Range("I113").Value = 2069.9
$ws.Range("K113").Value = 2069.9
$ws.Range("M113").Value = 100.0999999999999
$ws.Range("H126").Value = 1825
$ws.Range("I126").Value = 1828.5714
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 5485.7142
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -3015.7142
$ws.Range("N126").Value = -10340

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 455
$ws.Range("I107").Value = 552.3333
$ws.Range("J107").Value = 357.66666
$ws.Range("K107").Value = 1656.9999
$ws.Range("L107").Value = 1072.99998
$ws.Range("M107").Value = 263.0001
$ws.Range("N107").Value = -4912.999980000001
$ws.Range("H122").Value = 3359.3
$ws.Range("I122").Value = 1502
$ws.Range("K122").Value = 4506
$ws.Range("M122").Value = -2056
$ws.Range("H132").Value = 48822.74
$ws.Range("I132").Value = 34191.9
$ws.Range("J132").Value = 85399.836
$ws.Range("K132").Value = 102575.7
$ws.Range("L132").Value = 256199.508
$ws.Range("M132").Value = -100045.7
$ws.Range("N132").Value = -261259.508
